$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "ok" status markers for rows that now have a completed date (column E)
$ws.Range("E4").Value = "ok"
$ws.Range("E25").Value = "ok"
$ws.Range("E31").Value = "ok"
$ws.Range("E37").Value = "ok"
$ws.Range("E41").Value = "ok"

# New event added: "Sistema de Sumarios Administrativos" now has a link and missing-data note
$ws.Range("D38").Value = "falta data"
$ws.Range("E38").Value = "falta data!!!"

# Reset the view to the top of the sheet / select E1
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E1").Select() | Out-Null
